$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value2 = 1.999699333333333
$ws.Cells.Item(2, 8).Value2 = 5.999098
$ws.Cells.Item(2, 9).Value2 = 0.5605459117818491
$ws.Cells.Item(2, 10).Value2 = 0.5605459117818491
$ws.Cells.Item(2, 13).Value2 = 0.177232
$ws.Cells.Item(2, 14).Value2 = 0.531696
$ws.Cells.Item(2, 15).Value2 = 0.0005104719838156216
$ws.Cells.Item(2, 16).Value2 = 0.0005104719838156217
$ws.Cells.Item(2, 17).Value2 = 0.3544107122453333
$ws.Cells.Item(2, 18).Value2 = 3.189696410208
$ws.Cells.Item(2, 19).Value2 = 0.0002861429836070169
$ws.Cells.Item(2, 20).Value2 = 0.000286142983607017

$ws.Cells.Item(3, 7).Value2 = 1.999699333333333
$ws.Cells.Item(3, 8).Value2 = 5.999098
$ws.Cells.Item(3, 9).Value2 = 0.5605459117818491
$ws.Cells.Item(3, 10).Value2 = 0.5605459117818491
$ws.Cells.Item(3, 15).Value2 = 0.0002336907822601807
$ws.Cells.Item(3, 16).Value2 = 0.0002336907822601807
$ws.Cells.Item(3, 17).Value2 = 0.1622469385428889
$ws.Cells.Item(3, 18).Value2 = 1.460222446886
$ws.Cells.Item(3, 19).Value2 = 0.0001309944126170466
$ws.Cells.Item(3, 20).Value2 = 0.0001309944126170466

$ws.Cells.Item(4, 7).Value2 = 1.999699333333333
$ws.Cells.Item(4, 8).Value2 = 5.999098
$ws.Cells.Item(4, 9).Value2 = 0.5605459117818491
$ws.Cells.Item(4, 10).Value2 = 0.5605459117818491
$ws.Cells.Item(4, 13).Value2 = 274.5137023333334
$ws.Cells.Item(4, 14).Value2 = 823.541107
$ws.Cells.Item(4, 15).Value2 = 0.7906673411949746
$ws.Cells.Item(4, 16).Value2 = 0.7906673411949746
$ws.Cells.Item(4, 17).Value2 = 548.9448675468318
$ws.Cells.Item(4, 18).Value2 = 4940.503807921486
$ws.Cells.Item(4, 19).Value2 = 0.4432053456862674
$ws.Cells.Item(4, 20).Value2 = 0.4432053456862674

$ws.Cells.Item(5, 7).Value2 = 1.999699333333333
$ws.Cells.Item(5, 8).Value2 = 5.999098
$ws.Cells.Item(5, 9).Value2 = 0.5605459117818491
$ws.Cells.Item(5, 10).Value2 = 0.5605459117818491
$ws.Cells.Item(5, 13).Value2 = 0.042388
$ws.Cells.Item(5, 14).Value2 = 0.127164
$ws.Cells.Item(5, 15).Value2 = 0.0001220879211991998
$ws.Cells.Item(5, 16).Value2 = 0.0001220879211991998
$ws.Cells.Item(5, 17).Value2 = 0.08476325534133333
$ws.Cells.Item(5, 18).Value2 = 0.762869298072
$ws.Cells.Item(5, 19).Value2 = 0.00006843588510615597
$ws.Cells.Item(5, 20).Value2 = 0.00006843588510615597

$ws.Cells.Item(6, 7).Value2 = 1.999699333333333
$ws.Cells.Item(6, 8).Value2 = 5.999098
$ws.Cells.Item(6, 9).Value2 = 0.5605459117818491
$ws.Cells.Item(6, 10).Value2 = 0.5605459117818491
$ws.Cells.Item(6, 13).Value2 = 72.37795533333333
$ws.Cells.Item(6, 14).Value2 = 217.133866
$ws.Cells.Item(6, 15).Value2 = 0.2084664081177503
$ws.Cells.Item(6, 16).Value2 = 0.2084664081177503
$ws.Cells.Item(6, 17).Value2 = 144.7341490280964
$ws.Cells.Item(6, 18).Value2 = 1302.607341252868
$ws.Cells.Item(6, 19).Value2 = 0.1168549928142514
$ws.Cells.Item(6, 20).Value2 = 0.1168549928142514

$ws.Cells.Item(7, 7).Value2 = 1.567714666666667
$ws.Cells.Item(7, 8).Value2 = 4.703144
$ws.Cells.Item(7, 9).Value2 = 0.4394540882181509
$ws.Cells.Item(7, 10).Value2 = 0.4394540882181509
$ws.Cells.Item(7, 13).Value2 = 0.177232
$ws.Cells.Item(7, 14).Value2 = 0.531696
$ws.Cells.Item(7, 15).Value2 = 0.0005104719838156216
$ws.Cells.Item(7, 16).Value2 = 0.0005104719838156217
$ws.Cells.Item(7, 17).Value2 = 0.2778492058026666
$ws.Cells.Item(7, 18).Value2 = 2.500642852224
$ws.Cells.Item(7, 19).Value2 = 0.0002243290002086047
$ws.Cells.Item(7, 20).Value2 = 0.0002243290002086047

$ws.Cells.Item(8, 7).Value2 = 1.567714666666667
$ws.Cells.Item(8, 8).Value2 = 4.703144
$ws.Cells.Item(8, 9).Value2 = 0.4394540882181509
$ws.Cells.Item(8, 10).Value2 = 0.4394540882181509
$ws.Cells.Item(8, 15).Value2 = 0.0002336907822601807
$ws.Cells.Item(8, 16).Value2 = 0.0002336907822601807
$ws.Cells.Item(8, 17).Value2 = 0.1271975746231111
$ws.Cells.Item(8, 18).Value2 = 1.144778171608
$ws.Cells.Item(8, 19).Value2 = 0.0001026963696431342
$ws.Cells.Item(8, 20).Value2 = 0.0001026963696431342

$ws.Cells.Item(9, 7).Value2 = 1.567714666666667
$ws.Cells.Item(9, 8).Value2 = 4.703144
$ws.Cells.Item(9, 9).Value2 = 0.4394540882181509
$ws.Cells.Item(9, 10).Value2 = 0.4394540882181509
$ws.Cells.Item(9, 13).Value2 = 274.5137023333334
$ws.Cells.Item(9, 14).Value2 = 823.541107
$ws.Cells.Item(9, 15).Value2 = 0.7906673411949746
$ws.Cells.Item(9, 16).Value2 = 0.7906673411949746
$ws.Cells.Item(9, 17).Value2 = 430.3591573489342
$ws.Cells.Item(9, 18).Value2 = 3873.232416140408
$ws.Cells.Item(9, 19).Value2 = 0.3474619955087072
$ws.Cells.Item(9, 20).Value2 = 0.3474619955087072

$ws.Cells.Item(10, 7).Value2 = 1.567714666666667
$ws.Cells.Item(10, 8).Value2 = 4.703144
$ws.Cells.Item(10, 9).Value2 = 0.4394540882181509
$ws.Cells.Item(10, 10).Value2 = 0.4394540882181509
$ws.Cells.Item(10, 13).Value2 = 0.042388
$ws.Cells.Item(10, 14).Value2 = 0.127164
$ws.Cells.Item(10, 15).Value2 = 0.0001220879211991998
$ws.Cells.Item(10, 16).Value2 = 0.0001220879211991998
$ws.Cells.Item(10, 17).Value2 = 0.06645228929066667
$ws.Cells.Item(10, 18).Value2 = 0.598070603616
$ws.Cells.Item(10, 19).Value2 = 0.00005365203609304379
$ws.Cells.Item(10, 20).Value2 = 0.00005365203609304379

$ws.Cells.Item(11, 7).Value2 = 1.567714666666667
$ws.Cells.Item(11, 8).Value2 = 4.703144
$ws.Cells.Item(11, 9).Value2 = 0.4394540882181509
$ws.Cells.Item(11, 10).Value2 = 0.4394540882181509
$ws.Cells.Item(11, 13).Value2 = 72.37795533333333
$ws.Cells.Item(11, 14).Value2 = 217.133866
$ws.Cells.Item(11, 15).Value2 = 0.2084664081177503
$ws.Cells.Item(11, 16).Value2 = 0.2084664081177503
$ws.Cells.Item(11, 17).Value2 = 113.4679821194115
$ws.Cells.Item(11, 18).Value2 = 1021.211839074704
$ws.Cells.Item(11, 19).Value2 = 0.09161141530349891
$ws.Cells.Item(11, 20).Value2 = 0.09161141530349891
